# Commit: "Model fitting based on data preprocessing"
#
# 1) On the "Forecast Results" sheet, the raw/interim `traffic_volume`
#    column (D) held long inline-string blobs (concatenated historical
#    series used during preprocessing). After the model-fitting step the
#    column is collapsed down to a single numeric placeholder (0) for
#    every data row (rows 2-13).
# 2) On the "Metrics" sheet, the refit model produced new MAE / MAPE /
#    RMSE scores, replacing the previous ones in A2:C2.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Results")
$wsMetrics  = $wb.Worksheets.Item("Metrics")

# Column D ("traffic_volume") for rows 2 through 13 -> numeric 0
for ($row = 2; $row -le 13; $row++) {
    $wsForecast.Cells.Item($row, 4).Value = 0
}

# Updated metrics on the "Metrics" sheet
$wsMetrics.Range("A2").Value = 8.47162255564942
$wsMetrics.Range("B2").Value = 15.73131618321155
$wsMetrics.Range("C2").Value = 9.917842793772406
